# Update the division problems in the first (and only) table.
# Each populated row (1, 5, 9, 13, 17) has 5 cells; we address cells by
# (row, column) rather than by old text, since several of the old/new
# values repeat elsewhere in the table (e.g. "52÷2=", "98÷6=", "42÷6=").

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="50÷2="},
    @{Row=1;  Col=2; New="98÷6="},
    @{Row=1;  Col=3; New="90÷9="},
    @{Row=1;  Col=4; New="13÷2="},
    @{Row=1;  Col=5; New="13÷6="},

    @{Row=5;  Col=1; New="42÷6="},
    @{Row=5;  Col=2; New="86÷9="},
    @{Row=5;  Col=3; New="52÷2="},
    @{Row=5;  Col=4; New="54÷5="},
    @{Row=5;  Col=5; New="21÷3="},

    @{Row=9;  Col=1; New="67÷6="},
    @{Row=9;  Col=2; New="36÷7="},
    @{Row=9;  Col=3; New="25÷9="},
    @{Row=9;  Col=4; New="93÷8="},
    @{Row=9;  Col=5; New="33÷9="},

    @{Row=13; Col=1; New="27÷4="},
    @{Row=13; Col=2; New="45÷8="},
    @{Row=13; Col=3; New="88÷2="},
    @{Row=13; Col=4; New="99÷3="},
    @{Row=13; Col=5; New="42÷6="},

    @{Row=17; Col=1; New="74÷6="},
    @{Row=17; Col=2; New="72÷7="},
    @{Row=17; Col=3; New="10÷7="},
    @{Row=17; Col=4; New="35÷9="},
    @{Row=17; Col=5; New="38÷6="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}

Write-Host "Updated $($updates.Count) cells"
